$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column C
$ws.Range("C1").Value = "total xp"

# Row 2: base case, literal 0 (mirrors B2 which is a literal 0)
$ws.Range("C2").Value = 0

# Row 3: first running-total formula (references C2 literal)
$ws.Range("C3").Formula = "=B3+C2"

# Rows 4-201: running total "=B{r}+C{r-1}", entered in the same blocks
# that the existing B-column shared formulas use (B4:B67, B68:B131,
# B132:B195, B196:B201) so Excel groups the new C formulas the same way
$ws.Range("C4:C67").Formula = "=B4+C3"
$ws.Range("C68:C131").Formula = "=B68+C67"
$ws.Range("C132:C195").Formula = "=B132+C131"
$ws.Range("C196:C201").Formula = "=B196+C195"

# Note cell next to the level-19 XP (row 20) explaining the formula quirk
$ws.Range("J20").Value = "*didn't divide by 4"

# Restore the active selection to J21, as in the saved workbook
$ws.Range("J21").Select()
